$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5131
$ws.Range("I3").Value = 5349
$ws.Range("C4").Value = 1811
$ws.Range("E4").Value = 1967
$ws.Range("G4").Value = 1440
$ws.Range("I4").Value = 1227
$ws.Range("I5").Value = 497
$ws.Range("I6").Value = 5833
$ws.Range("C7").Value = 28354
$ws.Range("E7").Value = 25971
$ws.Range("G7").Value = 24662
$ws.Range("I7").Value = 18037

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 144
$ws.Range("I4").Value = 72
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 581
$ws.Range("I8").Value = 1089
$ws.Range("I10").Value = 129
$ws.Range("I11").Value = 270
$ws.Range("I15").Value = 208
$ws.Range("I18").Value = 131
$ws.Range("I20").Value = 435
$ws.Range("I23").Value = 177
$ws.Range("I29").Value = 1145
$ws.Range("I33").Value = 826
$ws.Range("I36").Value = 238
$ws.Range("I37").Value = 580
$ws.Range("I41").Value = 79
$ws.Range("I42").Value = 607
$ws.Range("I43").Value = 144
$ws.Range("I44").Value = 131
$ws.Range("I45").Value = 42
$ws.Range("I47").Value = 121
$ws.Range("I48").Value = 245
$ws.Range("I50").Value = 82
$ws.Range("I51").Value = 201
$ws.Range("I54").Value = 391
$ws.Range("C63").Value = 243
$ws.Range("E63").Value = 314
$ws.Range("I63").Value = 71
$ws.Range("G66").Value = 61
$ws.Range("I67").Value = 707
$ws.Range("I72").Value = 67
$ws.Range("I73").Value = 157
$ws.Range("I78").Value = 256
$ws.Range("I83").Value = 380
$ws.Range("I84").Value = 151
$ws.Range("I85").Value = 815
$ws.Range("I87").Value = 37
$ws.Range("I89").Value = 207
$ws.Range("I91").Value = 203
$ws.Range("I92").Value = 52
$ws.Range("I95").Value = 291
$ws.Range("I96").Value = 192
$ws.Range("I97").Value = 140
$ws.Range("I99").Value = 336
$ws.Range("I100").Value = 29
$ws.Range("C101").Value = 28354
$ws.Range("E101").Value = 25971
$ws.Range("G101").Value = 24662
$ws.Range("I101").Value = 18037

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 581

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 307
$ws.Range("I4").Value = 63
$ws.Range("I5").Value = 30
$ws.Range("I6").Value = 349
$ws.Range("I7").Value = 1089

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 117
$ws.Range("I4").Value = 24
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 208

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 40
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I6").Value = 136
$ws.Range("I7").Value = 435

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I2").Value = 50
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 341
$ws.Range("I3").Value = 396
$ws.Range("I5").Value = 42
$ws.Range("I6").Value = 307
$ws.Range("I7").Value = 1145

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 188
$ws.Range("I3").Value = 308
$ws.Range("I6").Value = 258
$ws.Range("I7").Value = 826

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 77
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 181
$ws.Range("I3").Value = 189
$ws.Range("I6").Value = 162
$ws.Range("I7").Value = 580

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 158
$ws.Range("I7").Value = 607

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I2").Value = 30
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 41
$ws.Range("I4").Value = 10
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 44
$ws.Range("I4").Value = 23
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 87
$ws.Range("I3").Value = 86
$ws.Range("I4").Value = 27
$ws.Range("I6").Value = 185
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("G4").Value = 7
$ws.Range("G7").Value = 61

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 257
$ws.Range("I7").Value = 707

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 60
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 134
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 380

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I3").Value = 50
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 223
$ws.Range("I3").Value = 326
$ws.Range("I7").Value = 815

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I2").Value = 5
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I3").Value = 49
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 203

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 109
$ws.Range("I7").Value = 291

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 55
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 121
$ws.Range("I7").Value = 336

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I3").Value = 6
$ws.Range("I6").Value = 29
